$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3356766.64
$ws.Range("C9").Value = 522931.99
$ws.Range("D9").Value = 3879698.63
$ws.Range("E9").Value = 13.47867553310449
$ws.Range("F9").Value = 86.52132446689552
$ws.Range("G9").Value = -49.46123062561033
$ws.Range("H9").Value = -39.38143433742655
$ws.Range("I9").Value = 33701
$ws.Range("J9").Value = 1419
$ws.Range("K9").Value = 35120
$ws.Range("L9").Value = 24250
$ws.Range("M9").Value = 159.987572371134
$ws.Range("N9").Value = 9.226774241356184
